# Add the new "ODI Batting Extra" worksheet as the last (4th) tab, matching
# the position used by the other sheets (after "ODI Bowling").
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "ODI Batting Extra"

# --- Header row (bold, centered, top-aligned, thin border -> matches the
#     header style used on the other sheets in this workbook) ---
$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$hdr = $ws.Range("A1:F1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108  # xlCenter
$hdr.VerticalAlignment = -4160    # xlTop
$hdr.Borders.LineStyle = 1
$hdr.Borders.Weight = 2

# --- Data rows ---
# Format the data range as text first so values such as "4488", "1", "0"
# and the percent strings are kept as literal text instead of being
# auto-coerced into numbers, matching the source data.
$data = $ws.Range("A2:F3")
$data.NumberFormat = "@"

$ws.Range("A2").Value = "4488"
$ws.Range("C2").Value = "1"
$ws.Range("D2").Value = "0"
$ws.Range("E2").Value = "3.18%"
$ws.Range("F2").Value = "NO"

$ws.Range("A3").Value = "4491"
$ws.Range("C3").Value = "2"
$ws.Range("D3").Value = "0"
$ws.Range("E3").Value = "14.40%"
$ws.Range("F3").Value = "NO"

# Restore General format / default style on the data range (keeps the
# values entered above as text, since the text was already committed).
$data.NumberFormat = "General"
$data.Style = "Normal"

# BATTING_POSITION is a genuine number column.
$ws.Cells.Item(2, 2).Value = 8
$ws.Cells.Item(3, 2).Value = 8
